# FA-1991 & FA 1992 Done
#
# Adds a new "Sheet5" (a report tab holding the same exDay/exMonth/exYear
# style data as Sheet3) to the end of the workbook and updates the
# view/selection state of the existing Sheet3 / Sheet4 tabs.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Sheet3")

# Create the new worksheet at the end of the workbook and name it Sheet5
$lastIndex = $wb.Worksheets.Count()
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$ws5.Name = "Sheet5"

# Bring over Sheet3's formatting (cell styles) and values for the A1:F2 block
$ws3.Range("A1:F2").Copy()
$ws5.Range("A1:F2").PasteSpecial(-4122)
$ws5.Range("A1:F2").PasteSpecial(-4163)

# Row 2 uses a taller row height than the sheet default, matching Sheet3
$ws5.Rows.Item(2).RowHeight = 16

# Sheet3 becomes the active tab with its whole data range selected
$ws3.Activate()
$ws3.Range("A1:F2").Select()

# Sheet5 is left with a lingering selection outside of the used range
$ws5.Range("C7").Select()

# Restore Sheet3 as the active sheet/tab (activeTab index 2)
$ws3.Activate()
